$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.573.27'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.930.91'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'198.15"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').Value = "'594.71"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').Value = '2.929.68'
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('D11').Value = "'0.439"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.43%  '
$ws.Range('D13').Value = '3.468.49'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = "'28.32"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '76.471.76'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = "'0.0000188"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '2.928.64'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').Value = "'13.42"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.75%  '
$ws.Range('D20').Value = "'8.66"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.60%  '
$ws.Range('D21').Value = "'372.61"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.95%  '
$ws.Range('D22').Value = "'4.29"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.38%  '
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').Value = "'72.10"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '3.089.47'
$ws.Range('E26').Value = '  +1.87%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = "'9.56"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').Value = "'0.0000107"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('D30').Value = "'0.998"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = "'8.25"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.07%  '
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D33').Value = "'497.97"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').Value = "'165.60"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  +19.18%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = "'0.391"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.34%  '
$ws.Range('D40').Value = "'19.95"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'0.110"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.16%  '
$ws.Range('D43').Value = "'179.11"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('E44').Value = '  -3.71%  '
$ws.Range('D45').Value = "'1.64"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = "'1.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.64%  '
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').Value = "'3.87"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('D50').Value = "'2.29"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.87%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = "'22.31"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.90%  '
